$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that currently holds account "004958578" (ASSAKO) is Excel row 22.
# We need to insert a brand-new row above it for account "004487140" (VALMIR, 2081.95),
# pushing "004958578" and everything below down by one row.
$targetRow = 22

$ws.Rows.Item($targetRow).Insert()

# Column A holds account numbers with leading zeros, so it must be stored
# as text (not auto-converted to a number) - match the existing cells above/below.
$ws.Cells.Item($targetRow, 1).Value = "'004487140"

$ws.Cells.Item($targetRow, 2).Value = "VALMIR"
$ws.Cells.Item($targetRow, 3).Value = 2081.95
